$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 11; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 23; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 30; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 36; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 54; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 61; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 93; I = 'ba'; J = 'Appreciation' }
    @{ Row = 98; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 102; I = '%'; J = 'Uninterpretable' }
    @{ Row = 109; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 112; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 125; I = 'ba'; J = 'Appreciation' }
    @{ Row = 134; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 140; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 141; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 143; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 155; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 156; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 166; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 168; I = 'ba'; J = 'Appreciation' }
    @{ Row = 177; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 199; I = '%'; J = 'Uninterpretable' }
    @{ Row = 220; I = '%'; J = 'Uninterpretable' }
    @{ Row = 224; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 225; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 235; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 245; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 253; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 259; I = 'qy'; J = 'Yes-No-Question' }
    @{ Row = 277; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 282; I = '%'; J = 'Uninterpretable' }
    @{ Row = 283; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 286; I = '%'; J = 'Uninterpretable' }
    @{ Row = 292; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 299; I = '%'; J = 'Uninterpretable' }
    @{ Row = 307; I = '%'; J = 'Uninterpretable' }
    @{ Row = 309; I = '%'; J = 'Uninterpretable' }
    @{ Row = 326; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 328; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 341; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 343; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 351; I = '%'; J = 'Uninterpretable' }
    @{ Row = 352; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 359; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 364; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 373; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 400; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 406; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 411; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 423; I = '%'; J = 'Uninterpretable' }
    @{ Row = 428; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 431; I = 'ba'; J = 'Appreciation' }
    @{ Row = 449; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 453; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 457; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 465; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 468; I = '%'; J = 'Uninterpretable' }
    @{ Row = 471; I = '%'; J = 'Uninterpretable' }
    @{ Row = 491; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 501; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 515; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 523; I = '%'; J = 'Uninterpretable' }
    @{ Row = 531; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 535; I = '%'; J = 'Uninterpretable' }
    @{ Row = 548; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 550; I = '%'; J = 'Uninterpretable' }
    @{ Row = 558; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 576; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 590; I = '%'; J = 'Uninterpretable' }
    @{ Row = 597; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 598; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 614; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 641; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 642; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 644; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 658; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 660; I = 'ba'; J = 'Appreciation' }
    @{ Row = 678; I = '%'; J = 'Uninterpretable' }
    @{ Row = 680; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 693; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 698; I = 'ba'; J = 'Appreciation' }
    @{ Row = 702; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 709; I = '%'; J = 'Uninterpretable' }
    @{ Row = 719; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 720; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 724; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 726; I = 'sv'; J = 'Statement-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}